$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New / updated data for rows 2-25 (columns A-T), reflecting the refreshed TPM-based
# NATMI ligand-receptor statistics and the addition of the "Neutrophils" sending cluster.
$data = @(
    @("Inflammatory-Mac", "Cdh1", "Itgb7", "ECs", 3, 1, 0.6295016666666666, 1.888505, 0.4100416847468479, 0.425450884376791, 2, 1, 1.5407, 3.0814, 0.03356114336260917, 0.02282654192720314, 0.9698732178333332, 5.819239306999999, 0.01376146776643476, 0.009711572450192475),
    @("Inflammatory-Mac", "Cdh1", "Itgb7", "FAPs", 3, 1, 0.6295016666666666, 1.888505, 0.4100416847468479, 0.425450884376791, 3, 1, 0.9053026666666666, 2.715908, 0.01972025221169104, 0.02011903285273785, 0.5698895375044444, 5.129005837539999, 0.008086125440514546, 0.008559660320003033),
    @("Inflammatory-Mac", "Cdh1", "Itgb7", "Inflammatory-Mac", 3, 1, 0.6295016666666666, 1.888505, 0.4100416847468479, 0.425450884376791, 3, 1, 22.301646, 66.904938, 0.4857978442449273, 0.4956215916122303, 14.03892332641, 126.35030993769, 0.1991973665005768, 0.2108626444676562),
    @("Inflammatory-Mac", "Cdh1", "Itgb7", "MuSCs", 3, 1, 0.6295016666666666, 1.888505, 0.4100416847468479, 0.425450884376791, 2, 1, 1.189092, 2.378184, 0.02590204912269207, 0.01761722489342626, 0.74853539582, 4.49121237492, 0.01062091986066427, 0.00749526391117302),
    @("Inflammatory-Mac", "Cdh1", "Itgb7", "Neutrophils", 3, 1, 0.6295016666666666, 1.888505, 0.4100416847468479, 0.425450884376791, 3, 1, 9.824249, 29.472747, 0.2140020958778281, 0.2183296213102325, 6.184381119248333, 55.65943007323499, 0.08774977993310112, 0.09288853047208831),
    @("Inflammatory-Mac", "Cdh1", "Itgb7", "Resolving-Mac", 3, 1, 0.6295016666666666, 1.888505, 0.4100416847468479, 0.425450884376791, 3, 1, 10.14626633333333, 30.438799, 0.2210166151802524, 0.2254859874041698, 6.387091567277221, 57.48382410549499, 0.09062602524555645, 0.09593321275567801),
    @("MuSCs", "Cdh1", "Itgb7", "ECs", 2, 1, 0.1668095, 0.333619, 0.1086555477667986, 0.0751591860201062, 2, 1, 1.5407, 3.0814, 0.03356114336260917, 0.02282654192720314, 0.25700339665, 1.0280135866, 0.003646604415744357, 0.001715624310902414),
    @("MuSCs", "Cdh1", "Itgb7", "FAPs", 2, 1, 0.1668095, 0.333619, 0.1086555477667986, 0.0751591860201062, 3, 1, 0.9053026666666666, 2.715908, 0.01972025221169104, 0.02011903285273785, 0.1510130851753333, 0.9060785110519999, 0.002142714806160711, 0.001512130132723552),
    @("MuSCs", "Cdh1", "Itgb7", "Inflammatory-Mac", 2, 1, 0.1668095, 0.333619, 0.1086555477667986, 0.0751591860201062, 3, 1, 22.301646, 66.904938, 0.4857978442449273, 0.4956215916122303, 3.720126418437, 22.320758510622, 0.05278463087036248, 0.03725051539956473),
    @("MuSCs", "Cdh1", "Itgb7", "MuSCs", 2, 1, 0.1668095, 0.333619, 0.1086555477667986, 0.0751591860201062, 2, 1, 1.189092, 2.378184, 0.02590204912269207, 0.01761722489342626, 0.198351841974, 0.793407367896, 0.002814401335708632, 0.00132409628292307),
    @("MuSCs", "Cdh1", "Itgb7", "Neutrophils", 2, 1, 0.1668095, 0.333619, 0.1086555477667986, 0.0751591860201062, 3, 1, 9.824249, 29.472747, 0.2140020958778281, 0.2183296213102325, 1.6387780635655, 9.832668381392999, 0.02325251495084837, 0.01640947662175511),
    @("MuSCs", "Cdh1", "Itgb7", "Resolving-Mac", 2, 1, 0.1668095, 0.333619, 0.1086555477667986, 0.0751591860201062, 3, 1, 10.14626633333333, 30.438799, 0.2210166151802524, 0.2254859874041698, 1.692493613930166, 10.154961683581, 0.02401468138797406, 0.01694734327223732),
    @("Neutrophils", "Cdh1", "Itgb7", "ECs", 3, 1, 0.5130273333333334, 1.539082, 0.3341732090958447, 0.34673130228853, 2, 1, 1.5407, 3.0814, 0.03356114336260917, 0.02282654192720314, 0.7904212124666667, 4.7425272748, 0.01121523497840882, 0.007914676609162875),
    @("Neutrophils", "Cdh1", "Itgb7", "FAPs", 3, 1, 0.5130273333333334, 1.539082, 0.3341732090958447, 0.34673130228853, 3, 1, 0.9053026666666666, 2.715908, 0.01972025221169104, 0.02011903285273785, 0.4644450129395556, 4.180005116456, 0.006589979965760224, 0.006975898461815515),
    @("Neutrophils", "Cdh1", "Itgb7", "Inflammatory-Mac", 3, 1, 0.5130273333333334, 1.539082, 0.3341732090958447, 0.34673130228853, 3, 1, 22.301646, 66.904938, 0.4857978442449273, 0.4956215916122303, 11.441353976324, 102.972185786916, 0.1623406245831707, 0.1718475199020226),
    @("Neutrophils", "Cdh1", "Itgb7", "MuSCs", 3, 1, 0.5130273333333334, 1.539082, 0.3341732090958447, 0.34673130228853, 2, 1, 1.189092, 2.378184, 0.02590204912269207, 0.01761722489342626, 0.6100366978480001, 3.660220187088, 0.00865577087748822, 0.006108443330007597),
    @("Neutrophils", "Cdh1", "Itgb7", "Neutrophils", 3, 1, 0.5130273333333334, 1.539082, 0.3341732090958447, 0.34673130228853, 3, 1, 9.824249, 29.472747, 0.2140020958778281, 0.2183296213102325, 5.040108266472667, 45.360974398254, 0.07151376713273047, 0.07570171392505852),
    @("Neutrophils", "Cdh1", "Itgb7", "Resolving-Mac", 3, 1, 0.5130273333333334, 1.539082, 0.3341732090958447, 0.34673130228853, 3, 1, 10.14626633333333, 30.438799, 0.2210166151802524, 0.2254859874041698, 5.205311960279777, 46.84780764251799, 0.07385783155828633, 0.07818305006046287),
    @("Resolving-Mac", "Cdh1", "Itgb7", "ECs", 2, 0.6666666666666666, 0.2258753333333333, 0.6776260000000001, 0.1471295583905087, 0.1526586273145729, 2, 1, 1.5407, 3.0814, 0.03356114336260917, 0.02282654192720314, 0.3480061260666667, 2.0880367564, 0.004937836202021239, 0.003484668556945376),
    @("Resolving-Mac", "Cdh1", "Itgb7", "FAPs", 2, 0.6666666666666666, 0.2258753333333333, 0.6776260000000001, 0.1471295583905087, 0.1526586273145729, 3, 1, 0.9053026666666666, 2.715908, 0.01972025221169104, 0.02011903285273785, 0.2044855416008889, 1.840369874408, 0.002901431999255554, 0.003071343938195756),
    @("Resolving-Mac", "Cdh1", "Itgb7", "Inflammatory-Mac", 2, 0.6666666666666666, 0.2258753333333333, 0.6776260000000001, 0.1471295583905087, 0.1526586273145729, 3, 1, 22.301646, 66.904938, 0.4857978442449273, 0.4956215916122303, 5.037391724132001, 45.33652551718801, 0.07147522229081726, 0.0756609118429869),
    @("Resolving-Mac", "Cdh1", "Itgb7", "MuSCs", 2, 0.6666666666666666, 0.2258753333333333, 0.6776260000000001, 0.1471295583905087, 0.1526586273145729, 2, 1, 1.189092, 2.378184, 0.02590204912269207, 0.01761722489342626, 0.268586551864, 1.611519311184, 0.003810957048830946, 0.002689421369322575),
    @("Resolving-Mac", "Cdh1", "Itgb7", "Neutrophils", 2, 0.6666666666666666, 0.2258753333333333, 0.6776260000000001, 0.1471295583905087, 0.1526586273145729, 3, 1, 9.824249, 29.472747, 0.2140020958778281, 0.2183296213102325, 2.219055517624667, 19.971499658622, 0.03148603386114815, 0.03332990029133061),
    @("Resolving-Mac", "Cdh1", "Itgb7", "Resolving-Mac", 2, 0.6666666666666666, 0.2258753333333333, 0.6776260000000001, 0.1471295583905087, 0.1526586273145729, 3, 1, 10.14626633333333, 30.438799, 0.2210166151802524, 0.2254859874041698, 2.291791290130444, 20.626121611174, 0.03251807698843553, 0.03442238131579163)
)

$startRow = 2
$r = $startRow
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}

# Ensure the sheet dimension reflects the expanded data range (A1:T25).
$ws.Range("A1:T25").Select() | Out-Null